# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns on Sheet1
# to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "57.829.69" },
    @{ Cell = "E2"; Value = "  -3.00%  " },
    @{ Cell = "D3"; Value = "2.285.45" },
    @{ Cell = "E3"; Value = "  -2.49%  " },
    @{ Cell = "E4"; Value = "  +0.03%  " },
    @{ Cell = "D5"; Value = "'529.92" },
    @{ Cell = "E5"; Value = "  -5.01%  " },
    @{ Cell = "D6"; Value = "'131.29" },
    @{ Cell = "E6"; Value = "  -0.69%  " },
    @{ Cell = "E7"; Value = "  +0.04%  " },
    @{ Cell = "D8"; Value = "'0.584" },
    @{ Cell = "E8"; Value = "  +1.07%  " },
    @{ Cell = "D9"; Value = "2.282.60" },
    @{ Cell = "E9"; Value = "  -2.57%  " },
    @{ Cell = "D10"; Value = "'0.0990" },
    @{ Cell = "E10"; Value = "  -4.75%  " },
    @{ Cell = "D11"; Value = "'5.45" },
    @{ Cell = "E11"; Value = "  -2.04%  " },
    @{ Cell = "E12"; Value = "  -0.14%  " },
    @{ Cell = "E13"; Value = "  -3.19%  " },
    @{ Cell = "D14"; Value = "'23.41" },
    @{ Cell = "E14"; Value = "  -2.28%  " },
    @{ Cell = "D15"; Value = "2.691.28" },
    @{ Cell = "E15"; Value = "  -2.60%  " },
    @{ Cell = "D16"; Value = "57.841.08" },
    @{ Cell = "E16"; Value = "  -2.93%  " },
    @{ Cell = "E17"; Value = "  -3.55%  " },
    @{ Cell = "D18"; Value = "2.272.18" },
    @{ Cell = "E18"; Value = "  -2.97%  " },
    @{ Cell = "D19"; Value = "'10.48" },
    @{ Cell = "E19"; Value = "  -4.38%  " },
    @{ Cell = "D20"; Value = "'4.17" },
    @{ Cell = "E20"; Value = "  -5.84%  " },
    @{ Cell = "D21"; Value = "'310.34" },
    @{ Cell = "E21"; Value = "  -2.55%  " },
    @{ Cell = "D22"; Value = "'6.38" },
    @{ Cell = "E22"; Value = "  -3.55%  " },
    @{ Cell = "E23"; Value = "  -0.07%  " },
    @{ Cell = "D24"; Value = "'62.30" },
    @{ Cell = "E24"; Value = "  -2.44%  " },
    @{ Cell = "E25"; Value = "  -1.95%  " },
    @{ Cell = "E26"; Value = "  +0.15%  " },
    @{ Cell = "D27"; Value = "'7.95" },
    @{ Cell = "E27"; Value = "  -4.53%  " },
    @{ Cell = "D28"; Value = "'1.26" },
    @{ Cell = "E28"; Value = "  -6.81%  " },
    @{ Cell = "D29"; Value = "'170.71" },
    @{ Cell = "E29"; Value = "  -0.28%  " },
    @{ Cell = "E30"; Value = "  -5.80%  " },
    @{ Cell = "D31"; Value = "0.0₃0715" },
    @{ Cell = "E31"; Value = "  -4.25%  " },
    @{ Cell = "E32"; Value = "  -4.01%  " },
    @{ Cell = "D33"; Value = "'1.04" },
    @{ Cell = "E33"; Value = "  -5.07%  " },
    @{ Cell = "D34"; Value = "'0.378" },
    @{ Cell = "E34"; Value = "  -5.14%  " },
    @{ Cell = "E35"; Value = "  +0.01%  " },
    @{ Cell = "D36"; Value = "'17.70" },
    @{ Cell = "E36"; Value = "  -1.71%  " },
    @{ Cell = "E37"; Value = "  -0.04%  " },
    @{ Cell = "D38"; Value = "'1.23" },
    @{ Cell = "E38"; Value = "  -5.76%  " },
    @{ Cell = "E39"; Value = "  -4.54%  " },
    @{ Cell = "D40"; Value = "'38.48" },
    @{ Cell = "E40"; Value = "  -0.20%  " },
    @{ Cell = "E41"; Value = "  -5.45%  " },
    @{ Cell = "D42"; Value = "'140.97" },
    @{ Cell = "E42"; Value = "  -2.11%  " },
    @{ Cell = "D43"; Value = "'284.01" },
    @{ Cell = "E43"; Value = "  -9.33%  " },
    @{ Cell = "E44"; Value = "  -2.26%  " },
    @{ Cell = "D45"; Value = "'0.0947" },
    @{ Cell = "E45"; Value = "  -1.21%  " },
    @{ Cell = "D46"; Value = "'0.0493" },
    @{ Cell = "E46"; Value = "  -2.19%  " },
    @{ Cell = "D47"; Value = "'0.551" },
    @{ Cell = "E47"; Value = "  -1.89%  " },
    @{ Cell = "D48"; Value = "'17.88" },
    @{ Cell = "E48"; Value = "  -5.35%  " },
    @{ Cell = "D49"; Value = "'0.0208" },
    @{ Cell = "E49"; Value = "  -3.69%  " },
    @{ Cell = "D50"; Value = "'10.89" },
    @{ Cell = "E50"; Value = "  -1.46%  " },
    @{ Cell = "E51"; Value = "  -0.53%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
